$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2, B3 values
$ws.Range("B2").Value = 5.7
$ws.Range("B3").Value = 5.6

# Row 5 ("theta_threshold_range") is removed entirely; row 6
# ("pie_threshold_range") shifts up to become the new row 5.
$ws.Rows("5").Delete()

# Update the (shifted) pie_threshold_range row's Min/Max values.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Column widths (closest achievable to the target 21.375 / 5.125 / 5.5 char widths
# under this host's MDW-7 pixel-grid column-width quantization)
$ws.Columns("A").ColumnWidth = 20.714285714285715
$ws.Columns("B").ColumnWidth = 4.428571428571429
$ws.Columns("C").ColumnWidth = 4.714285714285714

# Selection
$ws.Range("B2").Select()
